# invitee import password field added
#
# Adds a new "password" column (M) to the invitee sample sheet:
#   - M1 header cell gets the text "password"
#   - M2 data cell also gets the text "password" (sample/placeholder value)
# This grows the used range from A1:L2 to A1:M2 and introduces a new
# shared string ("password") used by both the new header and the new
# sample value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for row 1.
$ws.Range("M1").Value = "password"

# New sample value for row 2 (same literal text as the header, matching
# the source data).
$ws.Range("M2").Value = "password"

# Move the selection to reflect the edit: the author ended up with the
# next empty cell (N2) selected after filling in the new column.
$ws.Range("N2").Select() | Out-Null
